# "Politicas e template para registro de nao conformidades"
#
# This adds, to the quality-process workbook:
#   - Plan1 (sheet1): a new activity "Comunicar resultado da avaliação de
#     qualidade" right after "Acompanhar aderência entre esperado/Realizado".
#   - Plan2 (sheet2): two new bullet activities, "Registrar não
#     conformidades" and "Comunicar resultados da avaliação da qualidade",
#     right after "Identificar não conformidade".
# and leaves the workbook with Plan2 as the active sheet (matching the
# author's last on-screen selection while editing that sheet).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Plan1: insert the new "Comunicar resultado da avaliação de
# qualidade" activity as a new row 11 (pushing the Documentos block
# below it down by one row).
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Plan1")
$ws1.Rows.Item(11).Insert()
$ws1.Cells.Item(11, 2).Value = "Comunicar resultado da avaliação de qualidade"

# ------------------------------------------------------------------
# Plan2: insert two new bullets in the "Atividades" list, right after
# "Identificar não conformidade" (row 7) and before "Prover feedback
# para equipe e gerente" (old row 8). A blank separator row is kept
# above "Identificar não conformidade" just like the existing layout.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Plan2")
$ws2.Rows.Item(7).Insert()
$ws2.Rows.Item(9).Insert()
$ws2.Rows.Item(10).Insert()
$ws2.Cells.Item(9, 3).Value = "Registrar não conformidades"
$ws2.Cells.Item(10, 3).Value = "Comunicar resultados da avaliação da qualidade"

# Column C now holds longer text than before, so it was widened by
# hand (losing its "best fit" auto width) to keep things readable.
$ws2.Columns.Item(3).ColumnWidth = 169.5

# Page setup picked up while reviewing the sheet.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Final on-screen state: selection left on Plan1!B12, but Plan2 is the
# active tab with the cursor on C10 (where the new text was typed).
[void]$ws1.Select()
[void]$ws1.Range("B12").Select()

[void]$ws2.Select()
[void]$ws2.Range("C10").Select()
